$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 808, pushing the existing 808:919 block down to 810:919
$ws.Rows("808:809").Insert()

# New row 808 — Fecha 45131, Sin especificar / Pintón
$ws.Range("A808").Value = 4
$ws.Range("B808").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C808").Value = "Los Lagos"
$ws.Range("D808").Value = 45131
$ws.Range("E808").Value = 10
$ws.Range("F808").Value = "Fruta"
$ws.Range("G808").Value = 100108
$ws.Range("H808").Value = "Tropicales y subtropicales"
$ws.Range("I808").Value = 100108006
$ws.Range("J808").Value = "Plátano"
$ws.Range("K808").Value = "Sin especificar"
$ws.Range("L808").Value = "Pintón"
$ws.Range("M808").Value = 500
$ws.Range("N808").Value = 17000
$ws.Range("O808").Value = 17000
$ws.Range("P808").Value = 17000
$ws.Range("Q808").Value = "`$/caja 20 kilos"
$ws.Range("R808").Value = "Ecuador"
$ws.Range("S808").Value = 850
$ws.Range("T808").Value = 20

# New row 809 — Fecha 45131, Sin especificar / Primera Pintón
$ws.Range("A809").Value = 4
$ws.Range("B809").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C809").Value = "Los Lagos"
$ws.Range("D809").Value = 45131
$ws.Range("E809").Value = 10
$ws.Range("F809").Value = "Fruta"
$ws.Range("G809").Value = 100108
$ws.Range("H809").Value = "Tropicales y subtropicales"
$ws.Range("I809").Value = 100108006
$ws.Range("J809").Value = "Plátano"
$ws.Range("K809").Value = "Sin especificar"
$ws.Range("L809").Value = "Primera Pintón"
$ws.Range("M809").Value = 1000
$ws.Range("N809").Value = 18500
$ws.Range("O809").Value = 19000
$ws.Range("P809").Value = 18750
$ws.Range("Q809").Value = "`$/caja 20 kilos"
$ws.Range("R809").Value = "Ecuador"
$ws.Range("S809").Value = 938
$ws.Range("T809").Value = 20
